$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# D12:D14 held their BSE codes as text; convert them to real numbers.
$ws.Range("D12").Value = 500770
$ws.Range("D13").Value = 540777
$ws.Range("D14").Value = 532234

# New row 15 - Ipca Laboratories Limited
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = "IPCALAB"
$ws.Range("C15").Value = "Ipca Laboratories Limited"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "524494"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = 0.49
$ws.Range("F15").Value = 1179.45
$ws.Range("G15").Value = 942597
$ws.Range("H15").Value = "day"
$ws.Range("I15").Value = "12/06/2024 10:32:29"

# New row 16 - Berger Paints (I) Limited
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "BERGEPAINT"
$ws.Range("C16").Value = "Berger Paints (i) Limited"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "509480"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = 1.2
$ws.Range("F16").Value = 499
$ws.Range("G16").Value = 2520688
$ws.Range("H16").Value = "day"
$ws.Range("I16").Value = "12/06/2024 10:32:29"
